$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the authoritative diff.
$updates = @{
    "D2" = "69.972.71"
    "E2" = "  -1.00%  "
    "D3" = "3.544.18"
    "E3" = "  -1.15%  "
    "E4" = "  -0.13%  "
    "D5" = "611.76"
    "E5" = "  +4.42%  "
    "D6" = "185.31"
    "E6" = "  -0.51%  "
    "D7" = "0.623"
    "E7" = "  +0.31%  "
    "D8" = "0.999"
    "E8" = "  -0.12%  "
    "E9" = "  -0.31%  "
    "D10" = "0.648"
    "E10" = "  -0.40%  "
    "D11" = "53.64"
    "E11" = "  -1.10%  "
    "D12" = "0.0000310"
    "E12" = "  -4.17%  "
    "D13" = "9.47"
    "E13" = "  -0.79%  "
    "D14" = "4.103.94"
    "E14" = "  -1.07%  "
    "D15" = "623.83"
    "E15" = "  +9.81%  "
    "D16" = "69.968.70"
    "E16" = "  -0.92%  "
    "D17" = "12.68"
    "E17" = "  +2.07%  "
    "B18" = "Chainlink"
    "C18" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D18" = "18.88"
    "E18" = "  -3.51%  "
    "B19" = "WrappedEther"
    "C19" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D19" = "3.537.16"
    "E19" = "  -0.83%  "
    "E20" = "  -0.25%  "
    "D21" = "0.992"
    "E21" = "  -2.12%  "
    "D22" = "17.65"
    "E22" = "  -0.13%  "
    "D23" = "4.73"
    "E23" = "  +1.62%  "
    "D24" = "101.21"
    "E24" = "  +5.69%  "
    "E25" = "  +0.12%  "
    "D26" = "3.00"
    "E26" = "  +1.76%  "
    "D27" = "11.01"
    "E27" = "  -4.31%  "
    "D28" = "9.57"
    "E28" = "  +4.80%  "
    "D29" = "33.47"
    "E29" = "  +3.65%  "
    "D30" = "7.03"
    "E30" = "  -3.79%  "
    "D31" = "12.23"
    "E31" = "  -1.69%  "
    "E32" = "  -0.34%  "
    "D33" = "63.46"
    "E33" = "  -2.13%  "
    "E34" = "  +14.72%  "
    "D35" = "3.24"
    "E35" = "  -2.44%  "
    "D36" = "532.50"
    "E36" = "  -5.38%  "
    "B37" = "TheGraph"
    "C37" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "D37" = "0.401"
    "E37" = "  -3.92%  "
    "B38" = "Dai"
    "C38" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D38" = "1.00"
    "E38" = "  +0.01%  "
    "D39" = "37.14"
    "E39" = "  -1.62%  "
    "D40" = "0.0₃0781"
    "E40" = "  -2.58%  "
    "D41" = "3.541.43"
    "E41" = "  +4.74%  "
    "E42" = "  +5.63%  "
    "D43" = "0.137"
    "E43" = "  +1.76%  "
    "D44" = "0.0458"
    "E44" = "  +2.88%  "
    "D45" = "2.93"
    "E45" = "  -1.07%  "
    "D46" = "0.143"
    "E46" = "  +4.39%  "
    "D47" = "3.37"
    "E47" = "  -5.76%  "
    "D48" = "9.17"
    "E48" = "  -3.45%  "
    "E49" = "  +0.27%  "
    "E50" = "  -1.47%  "
    "D51" = "135.60"
    "E51" = "  -1.59%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "0.623", "69.972.71")
    # are not silently re-typed as numbers by the COM value-setter, then drop
    # the temporary Text number-format back to the sheet default style so no
    # stray style index is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
